$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$categories = @{
    101 = "بطاقات"
    102 = "معلومة"
    103 = "حسابات"
    104 = "مساعدة"
    105 = "قروض"
    106 = "قروض"
    107 = "مساعدة"
    108 = "مساعدة"
    109 = "بطاقات"
    110 = "عقاري"
    111 = "قروض"
    112 = "قروض"
    113 = "بطاقات"
    114 = "قروض"
    115 = "مساعدة"
    116 = "معلومة"
    117 = "معلومة"
    118 = "شهادات"
    119 = "قروض"
    120 = "مساعدة"
    121 = "مساعدة"
    122 = "الذهب"
    123 = "الذهب"
    124 = "قروض"
    125 = "حسابات"
    126 = "بطاقات"
    127 = "مساعدة"
    128 = "حسابات"
    129 = "قروض"
    130 = "مساعدة"
    131 = "مساعدة"
    132 = "معلومة"
    133 = "تحويلات"
    134 = "معلومة"
    135 = "قروض"
    136 = "قروض"
    137 = "قروض"
    138 = "قروض"
    139 = "معلومة"
    140 = "مساعدة"
    141 = "معلومة"
    142 = "مساعدة"
    143 = "مساعدة"
    144 = "قروض"
    145 = "بطاقات"
    146 = "قروض"
    147 = "بطاقات"
    148 = "قروض"
    149 = "معلومة"
    150 = "عقاري"
}

foreach ($row in $categories.Keys) {
    $ws.Range("D" + $row).Value = $categories[$row]
}

$ws.Range("F145").Select()
